$wb = $excel.ActiveWorkbook

$sheet6Urls = @(
    "https://gigazine.net/news/20191121-skin-integrated-wireless-interfaces/",
    "https://gigazine.net/news/20191108-world-first-penis-transplant/",
    "https://gigazine.net/news/20191023-nikon-small-world-2019-photomicrography/",
    "https://gigazine.net/news/20191020-skin-on-interfaces-case/",
    "https://gigazine.net/news/20191009-playstation-5-release-controller/",
    "https://gigazine.net/news/20191004-apple-virtual-keyboard-haptic-technology/",
    "https://gigazine.net/news/20190926-brain-pay-attention-filter/",
    "https://gigazine.net/news/20190906-apple-event-rumors/",
    "https://gigazine.net/news/20190813-us-navy-no-touchscreens/",
    "https://gigazine.net/news/20190810-barefoot-sole/",
    "https://gigazine.net/news/20190724-apple-iphone-11-leak/",
    "https://gigazine.net/news/20190723-asus-rog-phone-ii/",
    "https://gigazine.net/news/20190712-noise-boost-performance/",
    "https://gigazine.net/news/20190626-vr-mental-health/",
    "https://gigazine.net/news/20190531-sensor-packed-glove/",
    "https://gigazine.net/news/20190517-prevent-pregnancy-pullout-method/",
    "https://gigazine.net/news/20190428-ntt-supersonic-drone-chokaigi2019/",
    "https://gigazine.net/news/20190402-nintendo-labo-toy-con-vr-kit-wind-bird/",
    "https://gigazine.net/news/20190314-mouse-in-video-game-vr/",
    "https://gigazine.net/news/20190207-paper-versus-screens/"
)

$sheet7Urls = @(
    "https://gigazine.net/news/20191127-boston-dynamics-robot-police-dog/",
    "https://gigazine.net/news/20191125-nasa-bruie/",
    "https://gigazine.net/news/20191108-first-biometric-yubikey/",
    "https://gigazine.net/news/20191025-rats-drive-cars-relaxing/",
    "https://gigazine.net/news/20191021-apt41-teamviewer/",
    "https://gigazine.net/news/20191021-mercedes-benz-app-exposed-information/",
    "https://gigazine.net/news/20191011-iran-military-killer-robot/",
    "https://gigazine.net/news/20191004-google-levi-smart-jacket/",
    "https://gigazine.net/news/20190927-headline/",
    "https://gigazine.net/news/20190926-shotwell-lawsuit-patent/",
    "https://gigazine.net/news/20190913-anker-eufy-robovac-l70-hybrid/",
    "https://gigazine.net/news/20190902-franklin-expedition-ship-wreck/",
    "https://gigazine.net/news/20190825-how-to-find-hidden-camera/",
    "https://gigazine.net/news/20190817-fate-dominate-grail-war/",
    "https://gigazine.net/news/20190801-zoom-contact-lense-by-blink/",
    "https://gigazine.net/news/20190614-starsky-robotics-unmanned-truck/",
    "https://gigazine.net/news/20190531-canada-uses-anti-spam-law/",
    "https://gigazine.net/news/20190516-whatsapp-hacked-just-a-call/",
    "https://gigazine.net/news/20190425-hacker-break-tracking-apps/",
    "https://gigazine.net/news/20190402-disney-use-science-female-personality/"
)

$ws6 = $wb.Worksheets.Item("触覚")
for ($i = 0; $i -lt $sheet6Urls.Length; $i++) {
    $ws6.Cells.Item($i + 3, 1).Value = $sheet6Urls[$i]
}

$ws7 = $wb.Worksheets.Item("遠隔操作")
for ($i = 0; $i -lt $sheet7Urls.Length; $i++) {
    $ws7.Cells.Item($i + 3, 1).Value = $sheet7Urls[$i]
}
